$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 299
$ws.Range("I4").Value = 299
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 299
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -185
$ws.Range("H5").Value = 209.8125
$ws.Range("I5").Value = 177.6923
$ws.Range("J5").Value = 349
$ws.Range("K5").Value = 177.6923
$ws.Range("L5").Value = 349
$ws.Range("M5").Value = -62.69229999999999
$ws.Range("N5").Value = -579
$ws.Range("H8").Value = 335
$ws.Range("I8").Value = 52.5
$ws.Range("J8").Value = 900
$ws.Range("K8").Value = 157.5
$ws.Range("L8").Value = 2700
$ws.Range("M8").Value = -18.5
$ws.Range("N8").Value = -2978
$ws.Range("H9").Value = 125111.125
$ws.Range("I9").Value = 250072.25
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 250072.25
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = -249903.25
$ws.Range("N9").Value = -488
$ws.Range("H40").Value = 4499.6665
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -1575
$ws.Range("H132").Value = 2325.2856
$ws.Range("I132").Value = 1746.1666
$ws.Range("J132").Value = 5800
$ws.Range("K132").Value = 5238.4998
$ws.Range("L132").Value = 17400
$ws.Range("M132").Value = -2708.4998
$ws.Range("N132").Value = -22460
$ws.Range("H138").Value = 4826.8945
$ws.Range("I138").Value = 4197
$ws.Range("J138").Value = 4994.8667
$ws.Range("K138").Value = 12591
$ws.Range("L138").Value = 14984.6001
$ws.Range("M138").Value = -7451
$ws.Range("N138").Value = -25264.6001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1775
$ws.Range("I13").Value = 550
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 550
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -406
$ws.Range("N13").Value = -3288
$ws.Range("H43").Value = 20376.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20376.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20376.5
$ws.Range("N43").Value = -21002.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H86").Value = 6322.4
$ws.Range("I86").Value = 1704
$ws.Range("J86").Value = 13250
$ws.Range("K86").Value = 1704
$ws.Range("L86").Value = 13250
$ws.Range("M86").Value = -581
$ws.Range("N86").Value = -15496
$ws.Range("H89").Value = 6322.4
$ws.Range("I89").Value = 1704
$ws.Range("J89").Value = 13250
$ws.Range("K89").Value = 8520
$ws.Range("L89").Value = 66250
$ws.Range("M89").Value = -2904
$ws.Range("N89").Value = -77482
$ws.Range("H99").Value = 999
$ws.Range("I99").Value = 999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 499
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10001000
$ws.Range("I6").Value = 20000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 20000000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -19999887
$ws.Range("N6").Value = -2226
$ws.Range("H7").Value = 120
$ws.Range("I7").Value = 60
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 60
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = 53
$ws.Range("N7").Value = -526
$ws.Range("H10").Value = 577.5
$ws.Range("I10").Value = 577.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 577.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -438.5
$ws.Range("N10").ClearContents()
$ws.Range("H17").Value = 9
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9
$ws.Range("N17").Value = -357
$ws.Range("H22").Value = 389.4737
$ws.Range("I22").Value = 394.11765
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 394.11765
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -44.11765000000003
$ws.Range("N22").Value = -1050
$ws.Range("H25").Value = 7771
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 10906.5
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 10906.5
$ws.Range("M25").Value = -1326
$ws.Range("N25").Value = -11254.5
$ws.Range("H33").Value = 4015.5
$ws.Range("I33").Value = 4015.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4015.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3636.5
$ws.Range("N33").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 300
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -788
$ws.Range("N7").ClearContents()
$ws.Range("H11").Value = 1329.75
$ws.Range("I11").Value = 1106.3334
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 3319.0002
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = -3179.0002
$ws.Range("N11").Value = -6280
$ws.Range("H62").Value = 9999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 29997
$ws.Range("N62").Value = -31369
$ws.Range("H65").Value = 9999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 89991
$ws.Range("N65").Value = -96855
$ws.Range("H134").Value = 4400
$ws.Range("I134").Value = 4400
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13200
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8130

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 125.1
$ws.Range("I2").Value = 150.14285
$ws.Range("J2").Value = 66.666664
$ws.Range("K2").Value = 150.14285
$ws.Range("L2").Value = 66.666664
$ws.Range("M2").Value = -37.14285000000001
$ws.Range("H18").Value = 7000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7586

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3999.75
$ws.Range("I3").Value = 4999.6665
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 4999.6665
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -4887.6665
$ws.Range("N3").Value = -1224
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H13").Value = 1500
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1360
$ws.Range("H15").Value = 3999.75
$ws.Range("I15").Value = 4999.6665
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 4999.6665
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = -4829.6665
$ws.Range("N15").Value = -1340
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H93").Value = 884.58826
$ws.Range("I93").Value = 876.2727
$ws.Range("J93").Value = 899.8333
$ws.Range("K93").Value = 876.2727
$ws.Range("L93").Value = 899.8333
$ws.Range("M93").Value = 371.7273

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
